# Redis Sizing Input workbook update:
#  - Add a "TLS" header column between "Port" and "Password"
#  - Replace the sample data row with new demo values
#  - Add a new (currently empty) styled row further down the sheet
#  - Resize a couple of columns to fit the new content
#  - Move the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: insert "TLS" before "Password", shifting the two
#     trailing headers one column to the right (D->E, E->F). Read old
#     values into variables first so the shift doesn't clobber itself.
$oldD1 = $ws.Range("D1").Value2
$oldE1 = $ws.Range("E1").Value2

$ws.Range("F1").Value = $oldE1
$ws.Range("E1").Value = $oldD1
$ws.Range("D1").Value = "TLS"

# --- Data row: new sample values, and drop the old Password/User data
#     (no longer present after the layout change).
$ws.Range("A2").Value = "demo"
$ws.Range("B2").Value = "192.168.29.50"
$ws.Range("C2").Value = 6379
$ws.Range("D2:F2").ClearContents()

# --- New empty-but-styled row further down the sheet.
$ws.Rows.Item(7).RowHeight = 17
$font = $ws.Range("B7").Font
$font.Size = 13
$font.Name = "Helvetica Neue"
$font.Color = 3158322

# --- Column widths: B grows to fit the long IP string, E is brand new
#     and sized for the "User (ACL)"-style long values, the rest are
#     left as they were. (Integer inputs land closest to the bestFit
#     widths Excel itself would have computed here.)
$ws.Columns.Item(2).ColumnWidth = 60
$ws.Columns.Item(5).ColumnWidth = 44

# --- Selection moves off the edited range.
$ws.Range("E12").Select()
